$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The withdrawal row currently at row 2 (date 45294 / buy-date -693594 / 2000 / EUR)
# actually belongs at the bottom of the table (row 10); every other row shifts
# up by one to fill the gap. (Currency labels are unaffected in meaning --
# the EUR row stays EUR, the USD rows stay USD -- only their row position
# changes.)

# Capture the row-2 values before we overwrite anything.
$savedA = $ws.Range("A2").Value2
$savedB = $ws.Range("B2").Value2
$savedC = $ws.Range("C2").Value2
$savedD = $ws.Cells.Item(2, 4).Text

# Shift rows 3..10 up into rows 2..9.
for ($r = 3; $r -le 10; $r++) {
    $dest = $r - 1
    $ws.Range("A$dest").Value = $ws.Range("A$r").Value2
    $ws.Range("B$dest").Value = $ws.Range("B$r").Value2
    $ws.Range("C$dest").Value = $ws.Range("C$r").Value2
    $ws.Range("D$dest").Value = $ws.Cells.Item($r, 4).Text
}

# Put the original row-2 data into the now-vacated last row (10).
$ws.Range("A10").Value = $savedA
$ws.Range("B10").Value = $savedB
$ws.Range("C10").Value = $savedC
$ws.Range("D10").Value = $savedD
